$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2's value to 35
$ws.Range("A2").Value = 35

# Remove row 3 (A3 = 31) entirely so the used range shrinks back to A1:A2
$ws.Rows.Item(3).Delete()
